$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.042.48'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.68%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.910.95'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.61%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.16%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.8407'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +10.14%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.19'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.80%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.14%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3235'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +5.67%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '26.74'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +4.98%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07065'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +3.29%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08041'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.85%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7531'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.60%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.893.79'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.00%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.229'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.46%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.01'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.24%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.21'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.31%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.035.08'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.66%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.964'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.66%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '245.06'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.96%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007779'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.03%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.157.38'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.58%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9996'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.20%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9998'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.16%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.011'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.07%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1619'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +24.73%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '169.69'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.86%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.279'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.70%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.98'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.58%  '

$ws.Range('E29').Value = '  +2.91%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.375'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.28%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.518'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.11%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.310'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.35%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05654'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +7.82%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.102'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.37%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.287'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.72%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7368'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.55%  '

$ws.Range('E37').Value = '  +0.14%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01921'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.11%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.791'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.45%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4449'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.75%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '72.60'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.92%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.022'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.17%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8435'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.42%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.906'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.35%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9998'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.21%  '

$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '101.42'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.53%  '

$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.627'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.12%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.827'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.63%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '988.47'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +9.34%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.063.03'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.59%  '

$ws.Range('E51').Value = '  +0.95%  '
